# Daily attendance processing - 2026-01-25 06:45:57
# Swap the order of names in the "Recorded By" (column G) cells that list
# both the user email and "System" — change "dnasr281@gmail.com, System"
# to "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$startRow = $usedRange.Row

# Column G is the "Recorded By" column.
$col = 7

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
